# Quarterly database update + read_price algorithm change
# -----------------------------------------------------------
# The report shows a rolling 10-quarter window (columns E..N) for each
# metric. This edit rolls the window forward by one quarter:
#   - drops the oldest quarter ("فصل دوم منتهی به 1399/06")
#   - every remaining quarter's figures shift one column to the left
#   - the newest quarter ("فصل چهارم منتهی به 1401/12") is appended in
#     column N with freshly computed figures
# A couple of historical cells (column J, the "فصل چهارم منتهی به 1400/12"
# period) were also revised because of a read_price algorithm change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# --- Rolling quarter headers (row 8 and row 24) ---
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $headers[$i]
    $ws.Range($cols[$i] + "24").Value = $headers[$i]
}

# --- Helper to write a full E..N row of values ---
function Set-RowValues($rowNumber, $values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowNumber).Value = $values[$i]
    }
}

# هزینه حقوق و دستمزد
Set-RowValues 13 @(0, 0, 0, 0, 0, 2694, 0, 0, 0, 6927)

# هزینه مطالبات مشکوک الوصول
Set-RowValues 15 @(-395, 505, 214, 324, -34, 303, 172, 54, 424, -432)

# سایر هزینه ها
Set-RowValues 16 @(292, 781, 395, 394, 394, 720, 60, 993, 4792, 4816)

# جمع (subtotal line within the first block)
Set-RowValues 17 @(20637, 9292, 21500, 23732, 4076, 43464, 48375, 14521, 3444, 95383)

# second "جمع" style line
Set-RowValues 19 @(3582, 15377, 15640, 15401, 31790, -4608, 8281, 34108, 83625, -41984)

# جمع کل (grand total = sum of rows 10..19)
Set-RowValues 20 @(24116, 25955, 37749, 39851, 36226, 42573, 56888, 49676, 92285, 64710)

# تعداد پرسنل غیر تولیدی شرکت
Set-RowValues 26 @(58, 58, 61, 59, 59, 53, 58, 57, 58, 49)

# تعداد پرسنل تولیدی شرکت
Set-RowValues 27 @(349, 349, 341, 346, 346, 346, 329, 330, 320, 315)
